{"js": "// Update the date line and all 25 division-problem answers in the table.\n// Each mapping is a unique, unambiguous old->new text pair (verified no\n// overlaps among old/new values), so a straightforward search+replace for\n// each pair is safe regardless of processing order.\nconst replacements = [\n  [\"2023-09-29 Friday\", \"2023-09-30 Saturday\"],\n  [\"83\u00f79=9, 2\", \"45\u00f79=5, 0\"],\n  [\"46\u00f78=5, 6\", \"81\u00f76=13, 3\"],\n  [\"73\u00f77=10, 3\", \"34\u00f76=5, 4\"],\n  [\"81\u00f74=20, 1\", \"64\u00f77=9, 1\"],\n  [\"25\u00f77=3, 4\", \"43\u00f79=4, 7\"],\n  [\"62\u00f74=15, 2\", \"54\u00f74=13, 2\"],\n  [\"74\u00f74=18, 2\", \"34\u00f73=11, 1\"],\n  [\"97\u00f74=24, 1\", \"46\u00f73=15, 1\"],\n  [\"98\u00f73=32, 2\", \"51\u00f72=25, 1\"],\n  [\"90\u00f73=30, 0\", \"84\u00f79=9, 3\"],\n  [\"57\u00f76=9, 3\", \"63\u00f74=15, 3\"],\n  [\"26\u00f74=6, 2\", \"21\u00f77=3, 0\"],\n  [\"80\u00f76=13, 2\", \"60\u00f78=7, 4\"],\n  [\"78\u00f78=9, 6\", \"35\u00f72=17, 1\"],\n  [\"52\u00f79=5, 7\", \"44\u00f74=11, 0\"],\n  [\"94\u00f76=15, 4\", \"42\u00f79=4, 6\"],\n  [\"17\u00f74=4, 1\", \"34\u00f72=17, 0\"],\n  [\"33\u00f77=4, 5\", \"90\u00f79=10, 0\"],\n  [\"92\u00f78=11, 4\", \"94\u00f77=13, 3\"],\n  [\"60\u00f75=12, 0\", \"52\u00f74=13, 0\"],\n  [\"42\u00f78=5, 2\", \"14\u00f74=3, 2\"],\n  [\"43\u00f74=10, 3\", \"34\u00f79=3, 7\"],\n  [\"41\u00f77=5, 6\", \"69\u00f73=23, 0\"],\n  [\"41\u00f75=8, 1\", \"58\u00f79=6, 4\"],\n  [\"40\u00f77=5, 5\", \"91\u00f78=11, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and all 25 division-problem answers in the table.\n# Each mapping is a unique, unambiguous old->new text pair (verified no\n# overlaps among old/new values), so a straightforward Find/Replace for\n# each pair is safe regardless of processing order.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-09-29 Friday\", \"2023-09-30 Saturday\"),\n    @(\"83\u00f79=9, 2\", \"45\u00f79=5, 0\"),\n    @(\"46\u00f78=5, 6\", \"81\u00f76=13, 3\"),\n    @(\"73\u00f77=10, 3\", \"34\u00f76=5, 4\"),\n    @(\"81\u00f74=20, 1\", \"64\u00f77=9, 1\"),\n    @(\"25\u00f77=3, 4\", \"43\u00f79=4, 7\"),\n    @(\"62\u00f74=15, 2\", \"54\u00f74=13, 2\"),\n    @(\"74\u00f74=18, 2\", \"34\u00f73=11, 1\"),\n    @(\"97\u00f74=24, 1\", \"46\u00f73=15, 1\"),\n    @(\"98\u00f73=32, 2\", \"51\u00f72=25, 1\"),\n    @(\"90\u00f73=30, 0\", \"84\u00f79=9, 3\"),\n    @(\"57\u00f76=9, 3\", \"63\u00f74=15, 3\"),\n    @(\"26\u00f74=6, 2\", \"21\u00f77=3, 0\"),\n    @(\"80\u00f76=13, 2\", \"60\u00f78=7, 4\"),\n    @(\"78\u00f78=9, 6\", \"35\u00f72=17, 1\"),\n    @(\"52\u00f79=5, 7\", \"44\u00f74=11, 0\"),\n    @(\"94\u00f76=15, 4\", \"42\u00f79=4, 6\"),\n    @(\"17\u00f74=4, 1\", \"34\u00f72=17, 0\"),\n    @(\"33\u00f77=4, 5\", \"90\u00f79=10, 0\"),\n    @(\"92\u00f78=11, 4\", \"94\u00f77=13, 3\"),\n    @(\"60\u00f75=12, 0\", \"52\u00f74=13, 0\"),\n    @(\"42\u00f78=5, 2\", \"14\u00f74=3, 2\"),\n    @(\"43\u00f74=10, 3\", \"34\u00f79=3, 7\"),\n    @(\"41\u00f77=5, 6\", \"69\u00f73=23, 0\"),\n    @(\"41\u00f75=8, 1\", \"58\u00f79=6, 4\"),\n    @(\"40\u00f77=5, 5\", \"91\u00f78=11, 3\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
